# Corrections following third round of review
# Remove the "subgenus" column (column AS) from the Materials worksheet.
# This deletes both the "subgenus" header (row 1) and the "${subgenus}"
# template value (row 2), shifting all subsequent columns one to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

$ws.Range("AS:AS").Delete()
